$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the category label in A8: "Enrolled in  PhD program" -> "Active PhD student"
$ws.Range("A8").Value = "Active PhD student"

# Match the active selection recorded in the saved view state
$ws.Range("A12").Select()
